# Natmi following Dr Hou advice:
# Re-ran the NATMI LR-pair analysis for Ntf3-Ntrk2 including the "ECs"
# (endothelial cells) cluster, turning the previous 2x2 FAPs/sCs cluster
# matrix (rows 2-5) into a full 3x3 ECs/FAPs/sCs x ECs/FAPs/sCs matrix
# (rows 2-10) with updated statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntf3"
$ws.Range("C2").Value = "Ntrk2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.788130666666667
$ws.Range("H2").Value = 5.364392
$ws.Range("I2").Value = 0.08755714261138148
$ws.Range("J2").Value = 0.08755714261138148
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.132555
$ws.Range("N2").Value = 0.397665
$ws.Range("O2").Value = 0.002244436732931371
$ws.Range("P2").Value = 0.002244436732931371
$ws.Range("Q2").Value = 0.23702566052
$ws.Range("R2").Value = 2.133230944680001
$ws.Range("S2").Value = 0.0001965164671074952
$ws.Range("T2").Value = 0.0001965164671074952

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntf3"
$ws.Range("C3").Value = "Ntrk2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.788130666666667
$ws.Range("H3").Value = 5.364392
$ws.Range("I3").Value = 0.08755714261138148
$ws.Range("J3").Value = 0.08755714261138148
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 45.53127133333334
$ws.Range("N3").Value = 136.593814
$ws.Range("O3").Value = 0.7709408010078718
$ws.Range("P3").Value = 0.7709408010078719
$ws.Range("Q3").Value = 81.41586256345423
$ws.Range("R3").Value = 732.7427630710881
$ws.Range("S3").Value = 0.0675013736587789
$ws.Range("T3").Value = 0.06750137365877891

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ntf3"
$ws.Range("C4").Value = "Ntrk2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.788130666666667
$ws.Range("H4").Value = 5.364392
$ws.Range("I4").Value = 0.08755714261138148
$ws.Range("J4").Value = 0.08755714261138148
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 13.395535
$ws.Range("N4").Value = 40.186605
$ws.Range("O4").Value = 0.2268147622591968
$ws.Range("P4").Value = 0.2268147622591968
$ws.Range("Q4").Value = 23.95296692990667
$ws.Range("R4").Value = 215.57670236916
$ws.Range("S4").Value = 0.01985925248549508
$ws.Range("T4").Value = 0.01985925248549508

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntf3"
$ws.Range("C5").Value = "Ntrk2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 11.451921
$ws.Range("H5").Value = 34.355763
$ws.Range("I5").Value = 0.560751794520949
$ws.Range("J5").Value = 0.560751794520949
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.132555
$ws.Range("N5").Value = 0.397665
$ws.Range("O5").Value = 0.002244436732931371
$ws.Range("P5").Value = 0.002244436732931371
$ws.Range("Q5").Value = 1.518009388155
$ws.Range("R5").Value = 13.662084493395
$ws.Range("S5").Value = 0.001258571925680002
$ws.Range("T5").Value = 0.001258571925680002

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ntf3"
$ws.Range("C6").Value = "Ntrk2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 11.451921
$ws.Range("H6").Value = 34.355763
$ws.Range("I6").Value = 0.560751794520949
$ws.Range("J6").Value = 0.560751794520949
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 45.53127133333334
$ws.Range("N6").Value = 136.593814
$ws.Range("O6").Value = 0.7709408010078718
$ws.Range("P6").Value = 0.7709408010078719
$ws.Range("Q6").Value = 521.4205223388981
$ws.Range("R6").Value = 4692.784701050083
$ws.Range("S6").Value = 0.432306437634582
$ws.Range("T6").Value = 0.432306437634582

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ntf3"
$ws.Range("C7").Value = "Ntrk2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 11.451921
$ws.Range("H7").Value = 34.355763
$ws.Range("I7").Value = 0.560751794520949
$ws.Range("J7").Value = 0.560751794520949
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.395535
$ws.Range("N7").Value = 40.186605
$ws.Range("O7").Value = 0.2268147622591968
$ws.Range("P7").Value = 0.2268147622591968
$ws.Range("Q7").Value = 153.404608572735
$ws.Range("R7").Value = 1380.641477154615
$ws.Range("S7").Value = 0.127186784960687
$ws.Range("T7").Value = 0.127186784960687

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Ntf3"
$ws.Range("C8").Value = "Ntrk2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.182390333333333
$ws.Range("H8").Value = 21.547171
$ws.Range("I8").Value = 0.3516910628676694
$ws.Range("J8").Value = 0.3516910628676694
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.132555
$ws.Range("N8").Value = 0.397665
$ws.Range("O8").Value = 0.002244436732931371
$ws.Range("P8").Value = 0.002244436732931371
$ws.Range("Q8").Value = 0.952061750635
$ws.Range("R8").Value = 8.568555755715
$ws.Range("S8").Value = 0.0007893483401438733
$ws.Range("T8").Value = 0.0007893483401438733

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Ntf3"
$ws.Range("C9").Value = "Ntrk2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.182390333333333
$ws.Range("H9").Value = 21.547171
$ws.Range("I9").Value = 0.3516910628676694
$ws.Range("J9").Value = 0.3516910628676694
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 45.53127133333334
$ws.Range("N9").Value = 136.593814
$ws.Range("O9").Value = 0.7709408010078718
$ws.Range("P9").Value = 0.7709408010078719
$ws.Range("Q9").Value = 327.0233630889105
$ws.Range("R9").Value = 2943.210267800194
$ws.Range("S9").Value = 0.2711329897145109
$ws.Range("T9").Value = 0.2711329897145109

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Ntf3"
$ws.Range("C10").Value = "Ntrk2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.182390333333333
$ws.Range("H10").Value = 21.547171
$ws.Range("I10").Value = 0.3516910628676694
$ws.Range("J10").Value = 0.3516910628676694
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 13.395535
$ws.Range("N10").Value = 40.186605
$ws.Range("O10").Value = 0.2268147622591968
$ws.Range("P10").Value = 0.2268147622591968
$ws.Range("Q10").Value = 96.21196109382834
$ws.Range("R10").Value = 865.9076498444549
$ws.Range("S10").Value = 0.07976872481301467
$ws.Range("T10").Value = 0.07976872481301467
